# Adding in decision tree modeling (Ensemble Voting row + Retrained Models column)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 8: Ensemble Voting (entered first so its string lands before the D1 header
#     string in the shared-strings table, matching the author's save order) ---
$ws.Range("A8").Value = "Ensemble Voting"
$ws.Range("A3").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Font.Bold = $true

$ws.Range("B8").Value = 0.81
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Font.Bold = $true

$ws.Range("C8").Value = ""
$ws.Range("A6").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# --- New column D: header "Retrained Models Based on Hyperparamter Search V2 " ---
$ws.Range("D1").Value = "Retrained Models Based on Hyperparamter Search V2 "

# Give D1 a border matching the table's side borders (left/right thin) without top/bottom.
# Start from a cell that already carries the "apply fill + border1" combo (A3) so the
# resulting style keeps applyFill, then trim the top/bottom edges off.
$ws.Range("A3").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
$ws.Range("D1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$ws.Range("D1").Font.Bold = $false

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 30.417
$ws.Columns.Item(4).ColumnWidth = 44.417

# --- Selection moved to D15 (matches final author session state) ---
$ws.Range("D15").Select()

$excel.CutCopyMode = 0
Write-Host "edit complete"
